$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46, shifting rows 46:177 down to 47:178
$ws.Rows(46).Insert()

# Populate the newly inserted row 46 with data
$ws.Range("A46").Value = 7
$ws.Range("B46").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C46").Value = "Ñuble"
$ws.Range("D46").Value = 44497
$ws.Range("D46").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E46").Value = 16
$ws.Range("F46").Value = 100112008
$ws.Range("G46").Value = "Coliflor"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 200
$ws.Range("K46").Value = 700
$ws.Range("L46").Value = 800
$ws.Range("M46").Value = 750
$ws.Range("N46").Value = "$/unidad"
$ws.Range("O46").Value = "Región Metropolitana"
$ws.Range("P46").Value = 750
$ws.Range("Q46").Value = 1
$ws.Range("R46").Value = "Hortaliza"
